$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the table with new checker/respondent data (keeping header row intact).
$ws.Range("A2").Value = "Darryl"
$ws.Range("B2").Value = "0813136498045"
$ws.Range("C2").Value = "daryl2032@gmail.com"

$ws.Range("A3").Value = "reinaldi sianturi"
$ws.Range("B3").Value = "085653667887"
$ws.Range("C3").Value = "reinald020@gmail.com"

$ws.Range("A4").Value = "zhaky hanif s."
$ws.Range("B4").Value = "081287892654"
$ws.Range("C4").Value = "zhaky.hanif@yahoo.com"

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "'yoril baskara"
$ws.Range("B5").Value = "081312546455"
$ws.Range("C5").Value = "yoril.bass@yahoo.com"

$ws.Range("A6").Value = "vincentius albert"
$ws.Range("B6").Value = "085946556659"
$ws.Range("C6").Value = "albertvin9@gmail.com"

$ws.Range("C6").Select()
